$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new attendance-date column (E) mirroring the formatting/width of column D
$ws.Columns.Item(5).ColumnWidth = 12

$ws.Range("E1").NumberFormat = "@"
$ws.Range("E1").Value = "02/01/2555"
$ws.Range("E2").Value = "ขาดเรียน"
$ws.Range("E3").Value = "มาเรียน"
